# Update gh-pages to output generated at 456a3b4
# Increment the "想去人数" (F column) count by 1 for a handful of rows
# that are duplicated across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
foreach ($r in 8, 9, 16, 29, 35) {
    $cell = $ws1.Cells.Item($r, 6)  # column F
    $cell.Value = $cell.Value2 + 1
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($r in 10, 11, 18, 19, 29, 36) {
    $cell = $ws4.Cells.Item($r, 6)  # column F
    $cell.Value = $cell.Value2 + 1
}
